# Applies the "Updated cryptos list" data refresh described by the diff.
# All target cells are plain text (inlineStr) in the source workbook, including
# many that look numeric (e.g. "1.00", "0.197", "98.046.20"). Assigning such
# strings directly to Range.Value causes this engine (like real Excel) to
# auto-coerce them into numbers, which would change both the stored type and
# the displayed text (e.g. dropping trailing zeros). To keep them as literal
# text - matching the original workbook's cell typing - we temporarily switch
# the cell to the Text number format before writing the value, then restore
# the cell's original style so no formatting/style metadata changes leak into
# the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws 'D2' '98.046.20'
Set-TextValue $ws 'E2' '  -0.93%  '
Set-TextValue $ws 'D3' '3.431.11'
Set-TextValue $ws 'E3' '  +3.88%  '
Set-TextValue $ws 'D4' '1.00'
Set-TextValue $ws 'E4' '  -0.06%  '
Set-TextValue $ws 'D5' '257.25'
Set-TextValue $ws 'E5' '  +1.04%  '
Set-TextValue $ws 'D6' '658.41'
Set-TextValue $ws 'E6' '  +5.52%  '
Set-TextValue $ws 'D7' '1.49'
Set-TextValue $ws 'E7' '  +2.69%  '
Set-TextValue $ws 'E8' '  +5.87%  '
Set-TextValue $ws 'D9' '1.07'
Set-TextValue $ws 'E9' '  +10.24%  '
Set-TextValue $ws 'E10' '  -0.10%  '
Set-TextValue $ws 'D11' '3.431.16'
Set-TextValue $ws 'E11' '  +4.00%  '
Set-TextValue $ws 'E12' '  +6.80%  '
Set-TextValue $ws 'D13' '42.36'
Set-TextValue $ws 'E13' '  +6.33%  '
Set-TextValue $ws 'D14' '6.53'
Set-TextValue $ws 'E14' '  +19.54%  '
Set-TextValue $ws 'E15' '  +4.24%  '
Set-TextValue $ws 'D16' '97.781.50'
Set-TextValue $ws 'E16' '  -0.90%  '
Set-TextValue $ws 'D17' '4.074.41'
Set-TextValue $ws 'E17' '  +3.88%  '
Set-TextValue $ws 'D18' '8.81'
Set-TextValue $ws 'E18' '  +40.29%  '
Set-TextValue $ws 'D19' '3.427.31'
Set-TextValue $ws 'E19' '  +3.92%  '
Set-TextValue $ws 'D20' '17.75'
Set-TextValue $ws 'E20' '  +14.91%  '
Set-TextValue $ws 'D21' '0.518'
Set-TextValue $ws 'E21' '  +67.64%  '
Set-TextValue $ws 'E22' '  +16.82%  '
Set-TextValue $ws 'D23' '3.48'
Set-TextValue $ws 'E23' '  +0.50%  '
Set-TextValue $ws 'D24' '512.39'
Set-TextValue $ws 'E24' '  +5.71%  '
Set-TextValue $ws 'D25' '0.0000208'
Set-TextValue $ws 'E25' '  +2.75%  '
Set-TextValue $ws 'D26' '6.20'
Set-TextValue $ws 'E26' '  +10.30%  '
Set-TextValue $ws 'D27' '99.16'
Set-TextValue $ws 'E27' '  +11.49%  '
Set-TextValue $ws 'D28' '12.83'
Set-TextValue $ws 'E28' '  +7.59%  '
Set-TextValue $ws 'D29' '0.155'
Set-TextValue $ws 'E29' '  +13.75%  '
Set-TextValue $ws 'D30' '11.51'
Set-TextValue $ws 'E30' '  +12.69%  '
Set-TextValue $ws 'B31' 'Cronos'
Set-TextValue $ws 'C31' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws 'D31' '0.197'
Set-TextValue $ws 'E31' '  +4.92%  '
Set-TextValue $ws 'B32' 'Dai'
Set-TextValue $ws 'C32' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws 'D32' '1.00'
Set-TextValue $ws 'E32' '  +0.10%  '
Set-TextValue $ws 'E33' '  +0.71%  '
Set-TextValue $ws 'E34' '  +22.14%  '
Set-TextValue $ws 'D35' '30.04'
Set-TextValue $ws 'E35' '  +8.15%  '
Set-TextValue $ws 'D36' '2.21'
Set-TextValue $ws 'E36' '  +14.05%  '
Set-TextValue $ws 'D37' '7.89'
Set-TextValue $ws 'E37' '  +9.68%  '
Set-TextValue $ws 'E38' '  +6.82%  '
Set-TextValue $ws 'E39' '  +14.93%  '
Set-TextValue $ws 'D40' '517.91'
Set-TextValue $ws 'E40' '  +6.01%  '
Set-TextValue $ws 'D41' '24.73'
Set-TextValue $ws 'E41' '  -0.29%  '
Set-TextValue $ws 'D42' '0.859'
Set-TextValue $ws 'E42' '  +9.65%  '
Set-TextValue $ws 'B43' 'MantraDAO'
Set-TextValue $ws 'C43' 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue $ws 'D43' '3.68'
Set-TextValue $ws 'E43' '  +1.61%  '
Set-TextValue $ws 'B44' 'VeChain'
Set-TextValue $ws 'C44' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D44' '0.0420'
Set-TextValue $ws 'E44' '  +26.72%  '
Set-TextValue $ws 'E45' '  +7.01%  '
Set-TextValue $ws 'D46' '5.48'
Set-TextValue $ws 'E46' '  +16.55%  '
Set-TextValue $ws 'D47' '8.24'
Set-TextValue $ws 'E47' '  +12.54%  '
Set-TextValue $ws 'E48' '  +0.00%  '
Set-TextValue $ws 'E49' '  +18.15%  '
Set-TextValue $ws 'E50' '  +8.02%  '
Set-TextValue $ws 'D51' '51.13'
Set-TextValue $ws 'E51' '  +11.34%  '